$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(15, 8).Value = 5727.9575
$ws.Cells.Item(15, 9).Value = 5727.9575
$ws.Cells.Item(15, 11).Value = 17183.8725
$ws.Cells.Item(15, 13).Value = -17014.8725
$ws.Cells.Item(19, 8).Value = 306.16666
$ws.Cells.Item(19, 9).Value = 385.16666
$ws.Cells.Item(19, 10).Value = 266.66666
$ws.Cells.Item(19, 11).Value = 385.16666
$ws.Cells.Item(19, 12).Value = 266.66666
$ws.Cells.Item(19, 13).Value = -210.16666
$ws.Cells.Item(19, 14).Value = -616.66666
$ws.Cells.Item(33, 8).Value = 226.1842
$ws.Cells.Item(33, 9).Value = 179.29578
$ws.Cells.Item(33, 10).Value = 892
$ws.Cells.Item(33, 11).Value = 179.29578
$ws.Cells.Item(33, 12).Value = 892
$ws.Cells.Item(33, 13).Value = 49.70421999999999
$ws.Cells.Item(33, 14).Value = -1350
$ws.Cells.Item(138, 8).Value = 1427.86
$ws.Cells.Item(138, 9).Value = 668.125
$ws.Cells.Item(138, 10).Value = 1934.35
$ws.Cells.Item(138, 11).Value = 2004.375
$ws.Cells.Item(138, 12).Value = 5803.049999999999
$ws.Cells.Item(138, 13).Value = 3135.625
$ws.Cells.Item(138, 14).Value = -16083.05
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 2487.3635
$ws.Cells.Item(2, 9).Value = 2487.3635
$ws.Cells.Item(2, 11).Value = 2487.3635
$ws.Cells.Item(2, 13).Value = -2374.3635
$ws.Cells.Item(45, 8).Value = 2582.4119
$ws.Cells.Item(45, 9).Value = 2447.6667
$ws.Cells.Item(45, 11).Value = 2447.6667
$ws.Cells.Item(45, 13).Value = -2070.6667
$ws.Cells.Item(61, 8).Value = 1153.2
$ws.Cells.Item(61, 9).Value = 882.931
$ws.Cells.Item(61, 10).Value = 1865.7273
$ws.Cells.Item(61, 11).Value = 882.931
$ws.Cells.Item(61, 12).Value = 1865.7273
$ws.Cells.Item(61, 13).Value = -670.931
$ws.Cells.Item(61, 14).Value = -2289.7273
$ws.Cells.Item(74, 8).Value = 1568.4255
$ws.Cells.Item(74, 9).Value = 1329
$ws.Cells.Item(74, 11).Value = 1329
$ws.Cells.Item(74, 13).Value = -455
$ws.Cells.Item(77, 8).Value = 1568.4255
$ws.Cells.Item(77, 9).Value = 1329
$ws.Cells.Item(77, 11).Value = 6645
$ws.Cells.Item(77, 13).Value = -2277
$ws.Cells.Item(113, 8).Value = 45996.332
$ws.Cells.Item(113, 10).Value = 45996.332
$ws.Cells.Item(113, 12).Value = 45996.332
$ws.Cells.Item(113, 14).Value = -54674.332
$ws.Cells.Item(114, 8).Value = 45912
$ws.Cells.Item(114, 10).Value = 45912
$ws.Cells.Item(114, 12).Value = 45912
$ws.Cells.Item(114, 14).Value = -54590
$ws.Cells.Item(115, 8).Value = 20000
$ws.Cells.Item(115, 10).Value = 20000
$ws.Cells.Item(115, 12).Value = 20000
$ws.Cells.Item(115, 14).Value = -23134
$ws.Cells.Item(116, 8).Value = 2487.3635
$ws.Cells.Item(116, 9).Value = 2487.3635
$ws.Cells.Item(116, 11).Value = 2487.3635
$ws.Cells.Item(116, 13).Value = -193.3634999999999
$ws.Cells.Item(132, 8).Value = 2381.5862
$ws.Cells.Item(132, 9).Value = 1306.2941
$ws.Cells.Item(132, 10).Value = 3904.9167
$ws.Cells.Item(132, 11).Value = 3918.8823
$ws.Cells.Item(132, 12).Value = 11714.7501
$ws.Cells.Item(132, 13).Value = -1388.8823
$ws.Cells.Item(132, 14).Value = -16774.7501
$ws.Cells.Item(136, 8).Value = 1153.2
$ws.Cells.Item(136, 9).Value = 882.931
$ws.Cells.Item(136, 10).Value = 1865.7273
$ws.Cells.Item(136, 11).Value = 2648.793
$ws.Cells.Item(136, 12).Value = 5597.1819
$ws.Cells.Item(136, 13).Value = -98.79300000000012
$ws.Cells.Item(136, 14).Value = -10697.1819
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 2487.3635
$ws.Cells.Item(3, 9).Value = 2487.3635
$ws.Cells.Item(3, 11).Value = 2487.3635
$ws.Cells.Item(3, 13).Value = -2373.3635
$ws.Cells.Item(108, 8).Value = 0
$ws.Cells.Item(108, 10).Value = 0
$ws.Cells.Item(108, 12).Value = 0
$ws.Cells.Item(108, 14).ClearContents()
$ws.Cells.Item(114, 8).Value = 18735.666
$ws.Cells.Item(114, 9).Value = 8621
$ws.Cells.Item(114, 10).Value = 20000
$ws.Cells.Item(114, 11).Value = 8621
$ws.Cells.Item(114, 12).Value = 20000
$ws.Cells.Item(114, 13).Value = -4282
$ws.Cells.Item(114, 14).Value = -28678
$ws.Cells.Item(115, 8).Value = 20000
$ws.Cells.Item(115, 10).Value = 20000
$ws.Cells.Item(115, 12).Value = 20000
$ws.Cells.Item(115, 14).Value = -23134
$ws.Cells.Item(134, 8).Value = 3453
$ws.Cells.Item(134, 9).Value = 3691
$ws.Cells.Item(134, 10).Value = 3401.261
$ws.Cells.Item(134, 11).Value = 11073
$ws.Cells.Item(134, 12).Value = 10203.783
$ws.Cells.Item(134, 13).Value = -8538
$ws.Cells.Item(134, 14).Value = -15273.783
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 4343.9893
$ws.Cells.Item(31, 9).Value = 2363.7368
$ws.Cells.Item(31, 11).Value = 2363.7368
$ws.Cells.Item(31, 13).Value = -2068.7368
$ws.Cells.Item(34, 8).Value = 4343.9893
$ws.Cells.Item(34, 9).Value = 2363.7368
$ws.Cells.Item(34, 11).Value = 2363.7368
$ws.Cells.Item(34, 13).Value = -2161.7368
$ws.Cells.Item(122, 8).Value = 151000
$ws.Cells.Item(122, 9).Value = 601000
$ws.Cells.Item(122, 10).Value = 1000
$ws.Cells.Item(122, 11).Value = 1803000
$ws.Cells.Item(122, 12).Value = 3000
$ws.Cells.Item(122, 13).Value = -1800550
$ws.Cells.Item(122, 14).Value = -7900
$ws.Cells.Item(132, 8).Value = 42247.914
$ws.Cells.Item(132, 9).Value = 1136.1666
$ws.Cells.Item(132, 11).Value = 3408.4998
$ws.Cells.Item(132, 13).Value = -878.4998000000001
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(9, 8).Value = 76929720
$ws.Cells.Item(9, 10).Value = 7195.5835
$ws.Cells.Item(9, 12).Value = 21586.7505
$ws.Cells.Item(9, 14).Value = -22034.7505
$ws.Cells.Item(74, 8).Value = 10797
$ws.Cells.Item(74, 9).Value = 1950
$ws.Cells.Item(74, 10).Value = 13746
$ws.Cells.Item(74, 11).Value = 5850
$ws.Cells.Item(74, 12).Value = 41238
$ws.Cells.Item(74, 13).Value = -4789
$ws.Cells.Item(74, 14).Value = -43360
$ws.Cells.Item(77, 8).Value = 10797
$ws.Cells.Item(77, 9).Value = 1950
$ws.Cells.Item(77, 10).Value = 13746
$ws.Cells.Item(77, 11).Value = 17550
$ws.Cells.Item(77, 12).Value = 123714
$ws.Cells.Item(77, 13).Value = -12246
$ws.Cells.Item(77, 14).Value = -134322
$ws.Cells.Item(131, 8).Value = 4047.1943
$ws.Cells.Item(131, 9).Value = 12947.375
$ws.Cells.Item(131, 10).Value = 1504.2858
$ws.Cells.Item(131, 11).Value = 38842.125
$ws.Cells.Item(131, 12).Value = 4512.857400000001
$ws.Cells.Item(131, 13).Value = -33802.125
$ws.Cells.Item(131, 14).Value = -14592.8574
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(122, 8).Value = 977.7857
$ws.Cells.Item(122, 9).Value = 973.4167
$ws.Cells.Item(122, 10).Value = 1004
$ws.Cells.Item(122, 11).Value = 2920.2501
$ws.Cells.Item(122, 12).Value = 3012
$ws.Cells.Item(122, 13).Value = -470.2501000000002
$ws.Cells.Item(122, 14).Value = -7912
$ws.Cells.Item(126, 8).Value = 9813.308000000001
$ws.Cells.Item(126, 9).Value = 12377.3
$ws.Cells.Item(126, 10).Value = 1266.6666
$ws.Cells.Item(126, 11).Value = 37131.89999999999
$ws.Cells.Item(126, 12).Value = 3799.9998
$ws.Cells.Item(126, 13).Value = -34661.89999999999
$ws.Cells.Item(126, 14).Value = -8739.9998
$ws.Cells.Item(132, 8).Value = 3714.7368
$ws.Cells.Item(132, 9).Value = 2598.6667
$ws.Cells.Item(132, 11).Value = 7796.000100000001
$ws.Cells.Item(132, 13).Value = -5266.000100000001
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(61, 8).Value = 1537.8182
$ws.Cells.Item(61, 9).Value = 1654.0526
$ws.Cells.Item(61, 10).Value = 801.6667
$ws.Cells.Item(61, 11).Value = 1654.0526
$ws.Cells.Item(61, 12).Value = 801.6667
$ws.Cells.Item(61, 13).Value = -1452.0526
$ws.Cells.Item(61, 14).Value = -1205.6667
$ws.Cells.Item(113, 8).Value = 1537.8182
$ws.Cells.Item(113, 9).Value = 1654.0526
$ws.Cells.Item(113, 10).Value = 801.6667
$ws.Cells.Item(113, 11).Value = 1654.0526
$ws.Cells.Item(113, 12).Value = 801.6667
$ws.Cells.Item(113, 13).Value = 515.9474
$ws.Cells.Item(113, 14).Value = -5141.6667
$ws.Cells.Item(127, 8).Value = 0
$ws.Cells.Item(127, 10).Value = 0
$ws.Cells.Item(127, 12).Value = 0
$ws.Cells.Item(127, 14).ClearContents()
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(123, 8).Value = 49285
$ws.Cells.Item(123, 10).Value = 49285
$ws.Cells.Item(123, 12).Value = 49285
$ws.Cells.Item(123, 14).Value = -59085
